# Lanchkhuti Municipality area sheet: drop the "population census" sub-title
# row and the old 1989 / 2002 columns, keeping only the 2014 figure - i.e.
# restore the simpler single-year layout (see commit "fixed export and
# fixing maps").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 held the now-removed "(according to the population census data)"
# sub-title (shared string index 0 in the old file) - delete the whole row
# and let everything below shift up.
$ws.Rows.Item(2).Delete()

# Columns B and C held the 1989 and 2002 figures; only the 2014 column
# (old column D) should remain, becoming the new column B.
$ws.Columns.Item(2).Delete()
$ws.Columns.Item(2).Delete()

# The remaining rows (title / blank / "(sq. km)" / year header / Area row)
# now use a taller, uniform 20.1pt row height, and a few extra blank rows
# are kept below the table.
for ($r = 1; $r -le 8; $r++) {
    $ws.Rows.Item($r).RowHeight = 20.1
}
